$d = $word.ActiveDocument
$nbsp = [char]0x00A0

function Replace-Text($search, $replace) {
    $r = $d.Content.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $r) {
        Write-Host "WARNING: replace failed for: $search"
    }
    return $r
}

# --- Edit 1: "Any variable..." paragraph rewrite, and deletion of the
#     following "iterator" paragraph. ---

Replace-Text " Any variable which is reachable from the scope of the N1QL query can be referred to using " " Any variable in the scope of the function can be referred to using "

Replace-Text "`$" ":"

Replace-Text " syntax in the N1QL statement where parameters will be  substituted according to the rules of named parameters substitution in the N1QL grammar specification." " syntax in any part of the N1QL statement where parameters are allowed by N1QL in prepared queries."

# Delete the whole following paragraph ("The iterator we provide...local to the iterator.")
$iterParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -match "^The iterator we provide is an input iterator") {
        $iterParaIndex = $i
        break
    }
}
Write-Host "iterParaIndex=$iterParaIndex"
$delRange = $d.Paragraphs.Item($iterParaIndex).Range
Write-Host "delRange Start=$($delRange.Start) End=$($delRange.End)"
$delRange.Delete()

# --- Edit 2: Timers paragraph -- remove the page-break artifact (which
#     is achieved implicitly by any edit spanning the run boundary) and
#     move the "_GoBack" bookmark here (splitting "must" into "mu"/"st"). ---

$rngT = $d.Content
$foundT = $rngT.Find.Execute("must be smaller than system defined limits", $true, $false, $false, $false, $false, $true, 1, $false, "must be smaller than system defined limits", 2)
Write-Host "foundT: $foundT Start=$($rngT.Start) End=$($rngT.End)"

# Also remove the old "_GoBack" bookmark location (it sat in the
# "Any variable..." paragraph before, but that area has already been
# rewritten above so the original bookmark is gone implicitly). Now
# place a fresh "_GoBack" bookmark splitting "must" -> "mu" | "st".
$bmPos = $rngT.Start + 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
Write-Host "bookmark added at $bmPos"

# --- Edit 3: Drop the stray <w:lastRenderedPageBreak/> that precedes the
#     leading non-breaking spaces of the "phoneverify[meta.id] = ..."
#     code line. ---

$rngN = $d.Content
$foundN = $rngN.Find.Execute("phoneverify", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
Write-Host "foundN: $foundN Start=$($rngN.Start)"
$nbspTarget = $d.Range($rngN.Start - 4, $rngN.Start - 2)
$nbspSearch = "$nbsp$nbsp"
$rN2 = $nbspTarget.Find.Execute($nbspSearch, $true, $false, $false, $false, $false, $true, 1, $false, $nbspSearch, 2)
Write-Host "rN2: $rN2"

# --- Edit 4: Drop the <w:lastRenderedPageBreak/> before the
#     "State and Statelessness" heading. ---

$rS = $d.Content.Find.Execute("State and Statelessness", $true, $false, $false, $false, $false, $true, 1, $false, "State and Statelessness", 2)
Write-Host "rS: $rS"

